$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1447.9
$ws.Range("I112").Value = 1149.8
$ws.Range("J112").Value = 1547.2667
$ws.Range("K112").Value = 3449.4
$ws.Range("L112").Value = 4641.800099999999
$ws.Range("M112").Value = -2341.4
$ws.Range("N112").Value = -6857.800099999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 871.9394
$ws.Range("I137").Value = 858.4167
$ws.Range("J137").Value = 879.6667
$ws.Range("K137").Value = 2575.2501
$ws.Range("L137").Value = 2639.0001
$ws.Range("M137").Value = -25.2501000000002
$ws.Range("N137").Value = -7739.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6890.602
$ws.Range("I32").Value = 6958.2715
$ws.Range("J32").Value = 6721.4287
$ws.Range("K32").Value = 6958.2715
$ws.Range("L32").Value = 6721.4287
$ws.Range("M32").Value = -6671.2715
$ws.Range("N32").Value = -7295.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 933.65216
$ws.Range("I61").Value = 813.67645
$ws.Range("J61").Value = 1273.5834
$ws.Range("K61").Value = 813.67645
$ws.Range("L61").Value = 1273.5834
$ws.Range("M61").Value = -601.67645
$ws.Range("N61").Value = -1697.5834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 930.18866
$ws.Range("I74").Value = 864.1905
$ws.Range("J74").Value = 1182.1818
$ws.Range("K74").Value = 864.1905
$ws.Range("L74").Value = 1182.1818
$ws.Range("M74").Value = 9.809499999999957
$ws.Range("N74").Value = -2930.1818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 930.18866
$ws.Range("I77").Value = 864.1905
$ws.Range("J77").Value = 1182.1818
$ws.Range("K77").Value = 4320.9525
$ws.Range("L77").Value = 5910.909000000001
$ws.Range("M77").Value = 47.04749999999967
$ws.Range("N77").Value = -14646.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 933.65216
$ws.Range("I136").Value = 813.67645
$ws.Range("J136").Value = 1273.5834
$ws.Range("K136").Value = 2441.02935
$ws.Range("L136").Value = 3820.7502
$ws.Range("M136").Value = 108.9706499999998
$ws.Range("N136").Value = -8920.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17327.555
$ws.Range("I134").Value = 1245.6123
$ws.Range("J134").Value = 73614.36
$ws.Range("K134").Value = 3736.8369
$ws.Range("L134").Value = 220843.08
$ws.Range("M134").Value = -1201.8369
$ws.Range("N134").Value = -225913.08

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2380.8635
$ws.Range("I31").Value = 2138.5417
$ws.Range("J31").Value = 3027.0557
$ws.Range("K31").Value = 2138.5417
$ws.Range("L31").Value = 3027.0557
$ws.Range("M31").Value = -1843.5417
$ws.Range("N31").Value = -3617.0557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2380.8635
$ws.Range("I34").Value = 2138.5417
$ws.Range("J34").Value = 3027.0557
$ws.Range("K34").Value = 2138.5417
$ws.Range("L34").Value = 3027.0557
$ws.Range("M34").Value = -1936.5417
$ws.Range("N34").Value = -3431.0557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1216.3549
$ws.Range("I58").Value = 1357
$ws.Range("J58").Value = 960.63635
$ws.Range("K58").Value = 1357
$ws.Range("L58").Value = 960.63635
$ws.Range("M58").Value = -1154
$ws.Range("N58").Value = -1366.63635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1410.475
$ws.Range("I132").Value = 1053.4717
$ws.Range("K132").Value = 3160.4151
$ws.Range("M132").Value = -630.4151000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1249.2297
$ws.Range("I134").Value = 1188.9833
$ws.Range("J134").Value = 1507.4286
$ws.Range("K134").Value = 3566.949900000001
$ws.Range("L134").Value = 4522.2858
$ws.Range("M134").Value = -1031.949900000001
$ws.Range("N134").Value = -9592.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1216.3549
$ws.Range("I136").Value = 1357
$ws.Range("J136").Value = 960.63635
$ws.Range("K136").Value = 4071
$ws.Range("L136").Value = 2881.90905
$ws.Range("M136").Value = -1521
$ws.Range("N136").Value = -7981.90905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1146.2222
$ws.Range("I45").Value = 825
$ws.Range("K45").Value = 2475
$ws.Range("M45").Value = -1943

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 5347
$ws.Range("J74").Value = 5721.1113
$ws.Range("L74").Value = 17163.3339
$ws.Range("N74").Value = -19285.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 5347
$ws.Range("J77").Value = 5721.1113
$ws.Range("L77").Value = 51490.00169999999
$ws.Range("N77").Value = -62098.00169999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 406.71796
$ws.Range("I122").Value = 178.83333
$ws.Range("J122").Value = 602.0476
$ws.Range("K122").Value = 1609.49997
$ws.Range("L122").Value = 5418.4284
$ws.Range("M122").Value = 840.5000300000002
$ws.Range("N122").Value = -10318.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4285.95
$ws.Range("J129").Value = 4934.9375
$ws.Range("L129").Value = 14804.8125
$ws.Range("N129").Value = -24804.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 3250
$ws.Range("I130").Value = 1500
$ws.Range("K130").Value = 4500
$ws.Range("M130").Value = 520

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 55803.684
$ws.Range("I136").Value = 113061.11
$ws.Range("J136").Value = 4272
$ws.Range("K136").Value = 339183.33
$ws.Range("L136").Value = 12816
$ws.Range("M136").Value = -334083.33
$ws.Range("N136").Value = -23016

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 29900.703
$ws.Range("I137").Value = 1683.7368
$ws.Range("J137").Value = 59685.277
$ws.Range("K137").Value = 5051.2104
$ws.Range("L137").Value = 179055.831
$ws.Range("M137").Value = 48.78960000000006
$ws.Range("N137").Value = -189255.831

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 47206.184
$ws.Range("I139").Value = 57048.89
$ws.Range("J139").Value = 2914
$ws.Range("K139").Value = 171146.67
$ws.Range("L139").Value = 8742
$ws.Range("M139").Value = -166006.67
$ws.Range("N139").Value = -19022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 109691.18
$ws.Range("I140").Value = 168816.61
$ws.Range("J140").Value = 3265.4
$ws.Range("K140").Value = 506449.83
$ws.Range("L140").Value = 9796.200000000001
$ws.Range("M140").Value = -501269.83
$ws.Range("N140").Value = -20156.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4533.3335
$ws.Range("I70").Value = 4125
$ws.Range("K70").Value = 4125
$ws.Range("M70").Value = -3855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4533.3335
$ws.Range("I73").Value = 4125
$ws.Range("K73").Value = 4125
$ws.Range("M73").Value = -3189

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1937.5778
$ws.Range("I132").Value = 1889.8649
$ws.Range("K132").Value = 5669.5947
$ws.Range("M132").Value = -3139.5947

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2088.1555
$ws.Range("I136").Value = 1139.5135
$ws.Range("K136").Value = 3418.5405
$ws.Range("M136").Value = -868.5405000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 755.9245
$ws.Range("I132").Value = 710.08887
$ws.Range("K132").Value = 2130.26661
$ws.Range("M132").Value = 399.7333899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1287.2
$ws.Range("I136").Value = 1380.5
$ws.Range("J136").Value = 914
$ws.Range("K136").Value = 4141.5
$ws.Range("L136").Value = 2742
$ws.Range("M136").Value = -1591.5
$ws.Range("N136").Value = -7842
